# Add a new "Spain" worksheet (Zettler market test data), modeled on the
# existing "Italy" sheet, and fill in the Spain-specific values.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Duplicate the Italy sheet (keeps identical formatting/merges/column
# widths) and place the copy immediately after Italy - this becomes Spain.
$italy.Copy($null, $italy)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Spain"

# Make room for three extra "Attached Functionality" rows before the
# trailing "Wg" / "Attached Functionality" rows (13 -> 14,15,16 new,
# old 14/15 shift to 17/18).
$ws.Range("A14").EntireRow.Insert()
$ws.Range("A14").EntireRow.Insert()
$ws.Range("A14").EntireRow.Insert()

# Copy formatting (style) of the row above down into the new rows, then
# set their text.
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A13").Copy($ws.Range("A15"))
$ws.Range("A13").Copy($ws.Range("A16"))

$ws.Range("A14").Value = "PZ4DS(Dect/Fault)"
$ws.Range("A15").Value = "Three PZ4DS(Dect/Fault)"
$ws.Range("A16").Value = "Two PZ4DS(Dect/Fault)"

# Spain-specific market name / user story reference.
$ws.Range("B2").Value = "Spain Market"
$ws.Range("B4").Value = "NGC-3103/T2048"

# Resize the columns to fit the new (longer) content - matches the
# best-fit widths Excel computed for the new text.
$ws.Columns("A").ColumnWidth = 24.333333333333332
$ws.Columns("B").ColumnWidth = 14.333333333333334
$ws.Columns("C").ColumnWidth = 12.0
$ws.Columns("D").ColumnWidth = 13.833333333333334

# The wrapped "Constants"/"Input Value" cells now wrap onto two lines
# given the new column widths, matching Excel's automatic row growth.
$ws.Rows("3").RowHeight = 28.8
$ws.Rows("4").RowHeight = 28.8

# Restore Italy's selection (it is no longer the active tab) and leave
# Spain as the active sheet/tab with its own selection.
$italy.Activate()
$italy.Range("G15").Select()

$ws.Activate()
$ws.Range("C4").Select()
